$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Student name correction (typo fix) ---
$ws.Range("E2").Value2 = "Thomas Bentley"

# --- Row 6: Gravity (bonus) criteria - scored & commented ---
$ws.Rows.Item(6).RowHeight = 86.25
$ws.Range("F6").Value2 = 10
$ws.Range("G6").Value2 = "Uses actual gravity interactions between bodies to produce rotation, movement, and orbits. It is awesome to explore what interactions occur and how they change based on position and other variables."

# --- Row 16: Handle tunneling - scored & commented ---
$ws.Range("G16").Value2 = "Collision detection prevents tunneling (also it is visually evident that no tunneling occures)."

# --- Row 17: Octree broad phase collision detection - scored & commented ---
$ws.Range("F17").Value2 = 40
$ws.Range("G17").Value2 = "Octree is implemented(class Octree) and used(Source.cpp line $([char]0x2026).)"

# --- Row 20: Bezier surface scaling demonstration - commented ---
$ws.Range("G20").Value2 = "Bezier surface in main menu is scaled."

# --- Row 21: Bezier surface shearing demonstration - commented ---
$ws.Rows.Item(21).RowHeight = 29.25
$ws.Range("G21").Value2 = "Bezier surface is continuously sheared in a sinusoidal pattern."

# --- Row 52: Mouse interaction - commented ---
$ws.Range("G52").Value2 = "Menu is interacted with using the mouse. (Source.cpp line$([char]0x2026))"

# --- Row 54: Custom GUI system - commented ---
$ws.Rows.Item(54).RowHeight = 58.5
$ws.Range("G54").Value2 = "Menus were self implemented (not based off of an example or using an API), they are made up of buttons and images, and interacted with using the mouse."

# --- Row 55: Camera class via quaternions - scored & commented ---
$ws.Range("F55").Value2 = 10
$ws.Range("G55").Value2 = "The camera class uses SLERP to calculate the orientation of the camera (aka the rotation). This is possible because the way the camera is modeled can be thought of as a point inside of a circle and rotating to face different points of the outside of the circle. "

# --- Row 56: SLERP with quaternions - commented ---
$ws.Range("G56").Value2 = "Camera.cpp"

# --- Row 61: GLFW (not GLUT/freeglut) - commented ---
$ws.Range("G61").Value2 = "Source.cpp main method for example. There is no freeGLUT code used in this project."

# --- Row 68: Bonus justification - texturing comment replaces old "Textures" note ---
$ws.Range("G68").Value2 = "Used OpenGL texturing to texture some objects. Menus are not made of up simple colored objects, the images shown are made possible by texturing (with a texture shader, UVs, glGenTextures, etc), something we have not discussed in class."

# --- Update the active selection to reflect where the editor left off ---
$ws.Range("G16").Select()
